$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 2.68
$ws.Range("N2").Value = 4.4
$ws.Range("V2").Value = 1.59
$ws.Range("X2").Value = 16
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 19
$ws.Range("AA2").Value = 38
$ws.Range("AB2").Value = 13.5
$ws.Range("AC2").Value = 8.199999999999999
$ws.Range("AD2").Value = 12
$ws.Range("AE2").Value = 27
$ws.Range("AF2").Value = 20
$ws.Range("AG2").Value = 13
$ws.Range("AH2").Value = 16
$ws.Range("AJ2").Value = 44
$ws.Range("AK2").Value = 28
$ws.Range("AN2").Value = 22
$ws.Range("AO2").Value = 19
$ws.Range("F3").Value = 7.8
$ws.Range("K3").Value = 110
$ws.Range("S3").Value = 1.82
$ws.Range("V3").Value = 5.6
$ws.Range("H4").Value = 2.48
$ws.Range("Q4").Value = 1.73
$ws.Range("R4").Value = 1.43
$ws.Range("U4").Value = 2.3
$ws.Range("X4").Value = 1000
$ws.Range("F5").Value = 3.95
$ws.Range("G5").Value = 4.7
$ws.Range("H5").Value = 1.8
$ws.Range("J5").Value = 3.9
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 5
$ws.Range("T5").Value = 1.6
$ws.Range("U5").Value = 2.34
$ws.Range("W5").Value = 1.27
$ws.Range("X5").Value = 27
$ws.Range("Y5").Value = 14.5
$ws.Range("Z5").Value = 16
$ws.Range("AA5").Value = 23
$ws.Range("AB5").Value = 22
$ws.Range("AC5").Value = 11.5
$ws.Range("AD5").Value = 11.5
$ws.Range("AE5").Value = 21
$ws.Range("AF5").Value = 38
$ws.Range("AG5").Value = 21
$ws.Range("AK5").Value = 48
$ws.Range("AL5").Value = 48
$ws.Range("AM5").Value = 80
$ws.Range("AN5").Value = 38
$ws.Range("AO5").Value = 10
$ws.Range("H6").Value = 2.16
$ws.Range("I6").Value = 2.18
$ws.Range("N6").Value = 4.2
$ws.Range("O6").Value = 1.29
$ws.Range("P6").Value = 2.08
$ws.Range("Q6").Value = 1.89
$ws.Range("R6").Value = 1.41
$ws.Range("S6").Value = 3.3
$ws.Range("T6").Value = 1.75
$ws.Range("V6").Value = 1.84
$ws.Range("Y6").Value = 10.5
$ws.Range("AH6").Value = 17
$ws.Range("AL6").Value = 48
$ws.Range("AO6").Value = 14.5
$ws.Range("K7").Value = 3.55
$ws.Range("N7").Value = 2.94
$ws.Range("P7").Value = 1.66
$ws.Range("Q7").Value = 2.2
$ws.Range("R7").Value = 1.25
$ws.Range("S7").Value = 4.2
$ws.Range("T7").Value = 1.9
$ws.Range("U7").Value = 1.91
$ws.Range("Z7").Value = 17.5
$ws.Range("F8").Value = 2.44
$ws.Range("L8").Value = 1.4
$ws.Range("V8").Value = 1.47
$ws.Range("X8").Value = 1000
$ws.Range("F9").Value = 2.64
$ws.Range("G9").Value = 3.3
$ws.Range("H9").Value = 2.38
$ws.Range("I9").Value = 2.92
$ws.Range("J9").Value = 2.82
$ws.Range("K9").Value = 4.3
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 3.9
$ws.Range("P9").Value = 2.1
$ws.Range("Q9").Value = 1.6
$ws.Range("R9").Value = 1.47
$ws.Range("S9").Value = 2.5
$ws.Range("T9").Value = 1.58
$ws.Range("U9").Value = 2.34
$ws.Range("V9").Value = 1.56
$ws.Range("W9").Value = 1.43
$ws.Range("AC9").Value = 10.5
$ws.Range("F11").Value = 3.1
$ws.Range("G11").Value = 3.55
$ws.Range("I11").Value = 2.74
$ws.Range("AC11").Value = 1000
$ws.Range("L12").Value = 1.22
$ws.Range("H13").Value = 2.38
$ws.Range("Z13").Value = 19.5
$ws.Range("F14").Value = 1.79
$ws.Range("G14").Value = 1.87
$ws.Range("T14").Value = 1.76
$ws.Range("V14").Value = 1.25
$ws.Range("W14").Value = 2.14
$ws.Range("O15").Value = 1.31
$ws.Range("W15").Value = 1.62
$ws.Range("N16").Value = 3.55
$ws.Range("N17").Value = 1.1
$ws.Range("Q17").Value = 1.5
$ws.Range("S17").Value = 1.06
$ws.Range("H18").Value = 3.8
$ws.Range("I18").Value = 3.85
$ws.Range("S18").Value = 2.84
$ws.Range("U18").Value = 2.44
$ws.Range("V18").Value = 1.35
$ws.Range("Y18").Value = 17
$ws.Range("AA18").Value = 70
$ws.Range("AD18").Value = 15
$ws.Range("AE18").Value = 38
$ws.Range("I19").Value = 1.62
$ws.Range("J19").Value = 4.7
$ws.Range("K19").Value = 4.8
$ws.Range("Q19").Value = 1.58
$ws.Range("R19").Value = 1.67
$ws.Range("U19").Value = 2.42
$ws.Range("V19").Value = 2.6
$ws.Range("AC19").Value = 10.5
$ws.Range("AO19").Value = 6.2
$ws.Range("H20").Value = 1.96
$ws.Range("I20").Value = 1.98
$ws.Range("L20").Value = 1.31
$ws.Range("N20").Value = 5.6
$ws.Range("P20").Value = 2.56
$ws.Range("Q20").Value = 1.63
$ws.Range("R20").Value = 1.61
$ws.Range("S20").Value = 2.58
$ws.Range("V20").Value = 2.02
$ws.Range("AC20").Value = 9.4
$ws.Range("AO20").Value = 9
$ws.Range("F21").Value = 1.94
$ws.Range("G21").Value = 2.04
$ws.Range("H21").Value = 4.7
$ws.Range("I21").Value = 5.3
$ws.Range("K21").Value = 3.5
$ws.Range("M21").Value = 1.11
$ws.Range("N21").Value = 2.72
$ws.Range("O21").Value = 1.5
$ws.Range("P21").Value = 1.57
$ws.Range("Q21").Value = 2.48
$ws.Range("R21").Value = 1.21
$ws.Range("S21").Value = 5.1
$ws.Range("T21").Value = 2.16
$ws.Range("W21").Value = 1.96
$ws.Range("Y21").Value = 15
$ws.Range("AA21").Value = 160
$ws.Range("AG21").Value = 11
$ws.Range("H22").Value = 2.28
$ws.Range("I22").Value = 2.42
$ws.Range("J22").Value = 3.2
$ws.Range("N22").Value = 2.74
$ws.Range("O22").Value = 1.48
$ws.Range("R22").Value = 1.2
$ws.Range("S22").Value = 4.3
$ws.Range("U22").Value = 1.79
$ws.Range("V22").Value = 1.7
$ws.Range("X22").Value = 1000
$ws.Range("Y22").Value = 8.800000000000001
$ws.Range("AB22").Value = 12.5
$ws.Range("AC22").Value = 1000
$ws.Range("AD22").Value = 14.5
$ws.Range("G23").Value = 2.38
$ws.Range("H23").Value = 3.6
$ws.Range("P23").Value = 1.53
$ws.Range("Q23").Value = 2.4
$ws.Range("S23").Value = 1.05
$ws.Range("V23").Value = 1.32
$ws.Range("W23").Value = 1.72
$ws.Range("AB23").Value = 1000
